# Edit script: "Arreglos al unir certificados"
# Refreshes the CC/NT certificate roster data (rows 2-21), the sheet view
# (zoom + selection) and the page setup (landscape), matching the
# author's re-saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 93343273
$ws.Range("C2").Value = "NT"
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = "FEBRERO"
$ws.Range("F2").Value = 1983
$ws.Range("B3").Value = 1106309469
$ws.Range("C3").Value = "NT"
$ws.Range("D3").Value = 19
$ws.Range("E3").Value = "ABRIL"
$ws.Range("F3").Value = 2016
$ws.Range("B4").Value = 1109494165
$ws.Range("C4").Value = "NT"
$ws.Range("D4").Value = 17
$ws.Range("E4").Value = "NOVIEMBRE"
$ws.Range("F4").Value = 2011
$ws.Range("B5").Value = 28890562
$ws.Range("C5").Value = "NT"
$ws.Range("D5").Value = 30
$ws.Range("E5").Value = "JUNIO"
$ws.Range("F5").Value = 1971
$ws.Range("B6").Value = 1110513628
$ws.Range("C6").Value = "NT"
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = "AGOSTO"
$ws.Range("F6").Value = 2009
$ws.Range("B7").Value = 93154117
$ws.Range("C7").Value = "NT"
$ws.Range("D7").Value = 15
$ws.Range("E7").Value = "DICIEMBRE"
$ws.Range("F7").Value = 1998
$ws.Range("B8").Value = 28898397
$ws.Range("C8").Value = "NT"
$ws.Range("D8").Value = 26
$ws.Range("E8").Value = "SEPTIEMBRE"
$ws.Range("F8").Value = 1972
$ws.Range("B9").Value = 1005824385
$ws.Range("C9").Value = "NT"
$ws.Range("D9").Value = 28
$ws.Range("E9").Value = "FEBRERO"
$ws.Range("F9").Value = 2019
$ws.Range("B10").Value = 1109496000
$ws.Range("C10").Value = "NT"
$ws.Range("D10").Value = 26
$ws.Range("E10").Value = "ABRIL"
$ws.Range("F10").Value = 2016
$ws.Range("B11").Value = 1109496462
$ws.Range("C11").Value = "NT"
$ws.Range("D11").Value = 23
$ws.Range("E11").Value = "MAYO"
$ws.Range("F11").Value = 2017
$ws.Range("B12").Value = 1119583112
$ws.Range("C12").Value = "NT"
$ws.Range("D12").Value = 16
$ws.Range("E12").Value = "JUNIO"
$ws.Range("F12").Value = 2020
$ws.Range("B13").Value = 65587237
$ws.Range("C13").Value = "NT"
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = "NOVIEMBRE"
$ws.Range("F13").Value = 1991
$ws.Range("B14").Value = 65798919
$ws.Range("C14").Value = "NT"
$ws.Range("D14").Value = 7
$ws.Range("E14").Value = "FEBRERO"
$ws.Range("F14").Value = 1994
$ws.Range("B15").Value = 28852225
$ws.Range("C15").Value = "NT"
$ws.Range("D15").Value = 26
$ws.Range("E15").Value = "OCTUBRE"
$ws.Range("F15").Value = 1977
$ws.Range("B16").Value = 1007684243
$ws.Range("C16").Value = "NT"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = "ABRIL"
$ws.Range("F16").Value = 2019
$ws.Range("B17").Value = 65720516
$ws.Range("C17").Value = "NT"
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = "OCTUBRE"
$ws.Range("F17").Value = 1981
$ws.Range("B18").Value = 65586697
$ws.Range("C18").Value = "NT"
$ws.Range("D18").Value = 28
$ws.Range("E18").Value = "FEBRERO"
$ws.Range("F18").Value = 1989
$ws.Range("B19").Value = 1109841349
$ws.Range("C19").Value = "NT"
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = "JULIO"
$ws.Range("F19").Value = 2005
$ws.Range("B20").Value = 36302746
$ws.Range("C20").Value = "NT"
$ws.Range("D20").Value = 17
$ws.Range("E20").Value = "DICIEMBRE"
$ws.Range("F20").Value = 1999
$ws.Range("B21").Value = 65588663
$ws.Range("C21").Value = "NT"
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = "ABRIL"
$ws.Range("F21").Value = 1999

# Sheet view: zoomed to 90% with a plain cell selection (no leftover
# "topLeftCell" scroll anchor from the previous view).
[void]$ws.Range("F14").Select()
$excel.ActiveWindow.Zoom = 90

# Page setup: landscape orientation, paper size 257.
$ws.PageSetup.Orientation = 2
$ws.PageSetup.PaperSize = 257
